# "Error Calculations and Plots"
#
# The source data table (Sheet1, A1:F35) loses two sample rows entirely
# ("RM 232" and "SC 92") and several of the remaining rows' column-F
# ("F" = error) values are refreshed -- some newly populated, some newly
# blanked out. After the two row deletions the table shrinks from
# A1:F35 to A1:F33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RM 232" row (originally row 26). Everything below shifts up.
$ws.Rows.Item(26).Delete()

# Remove the "SC 92" row. It was originally row 28, but after the delete
# above it has already shifted up to row 27.
$ws.Rows.Item(27).Delete()

# Refresh column F ("error") values on the final, post-delete row numbers.
$ws.Range("F19").Value = 17.81   # RM 125: was blank, now populated
$ws.Range("F21").Value = ""      # RM 135: was 16.58, now blank
$ws.Range("F23").Value = 16.48   # RM 140: was blank, now populated
$ws.Range("F27").Value = ""      # SC 101: was 17, now blank
$ws.Range("F33").Value = 17.53   # SC 232: was blank, now populated
